# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 8 (pushing the existing
# rows 8..72 down to 9..73, and growing the used range to A1:R73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 8; Excel shifts every
# row below it down by one (old row 72 becomes the new row 73) and
# inherits the formatting (date style, etc.) of the row above.
$ws.Rows("8:8").Insert()

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44761
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112012
$ws.Range("G8").Value = "Espinaca"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 1750
$ws.Range("N8").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 583
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = "Hortaliza"
